# Update cryptocurrency price/volume data (Sheet1, columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.988.84"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.843.98"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'360.11"
$ws.Range("E5").Value = "  +5.79%  "
$ws.Range("D6").Value = "'113.13"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  +4.17%  "
$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "'41.51"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.0863"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "'7.85"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").Value = "3.294.05"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "2.846.33"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "'0.929"
$ws.Range("E17").Value = "  +5.98%  "
$ws.Range("D18").Value = "51.933.47"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  +8.42%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'13.56"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("D22").Value = "0.0₃0997"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("D23").Value = "'70.42"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").Value = "'269.67"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'2.85"
$ws.Range("E25").Value = "  +3.52%  "
$ws.Range("D26").Value = "'27.17"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D28").Value = "'10.39"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").Value = "'53.33"
$ws.Range("E30").Value = "  +6.49%  "
$ws.Range("D31").Value = "'35.71"
$ws.Range("E31").Value = "  +2.97%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "'0.0475"
$ws.Range("E33").Value = "  +25.17%  "
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("D35").Value = "'5.51"
$ws.Range("E35").Value = "  +11.52%  "
$ws.Range("D36").Value = "'0.0853"
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "'18.66"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "'23.47"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").Value = "'125.26"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("D46").Value = "'3.42"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "2.110.54"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D49").Value = "'6.02"
$ws.Range("E49").Value = "  +8.43%  "
$ws.Range("D50").Value = "'0.982"
$ws.Range("E50").Value = "  +11.48%  "
$ws.Range("D51").Value = "'62.15"
$ws.Range("E51").Value = "  +4.86%  "
